# Contract & Order - Update/Delete scenario
$wb = $excel.ActiveWorkbook

# Sheet: Contract - new columns for the "update account name" test data
# (ContractTermUpdated / UpdateAccountName), used by the update scenario.
$wsContract = $wb.Worksheets.Item("Contract")
$wsContract.Range("C1").Value = "ContractTermUpdated"
$wsContract.Range("D1").Value = "UpdateAccountName"
$wsContract.Range("C2").Value = 2
$wsContract.Range("D2").Value = "Kevin Testing"

# Sheet: Order - new Street column used by the order update scenario
$wsOrder = $wb.Worksheets.Item("Order")
$wsOrder.Range("E1").Value = "Street"
$wsOrder.Range("E2").Value = "Test brook street"
$wsOrder.Activate()
$wsOrder.Range("F8").Select() | Out-Null

# Sheet: Opportunity - leave data as-is, just move the selection
$wsOpportunity = $wb.Worksheets.Item("Opportunity")
$wsOpportunity.Activate()
$wsOpportunity.Range("B8").Select() | Out-Null

# Sheet: Quote - no longer the active tab, move the selection
$wsQuote = $wb.Worksheets.Item("Quote")
$wsQuote.Activate()
$wsQuote.Range("D15").Select() | Out-Null

# Sheet: LoginPage - rotate the automation password and become the active tab
$wsLogin = $wb.Worksheets.Item("LoginPage")
$wsLogin.Range("C2").Value = "Automation@March24"
$wsLogin.Activate()
$wsLogin.Range("D6").Select() | Out-Null
